# Regenerate advocacy packet figures (2026-02-22)

$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "Generated: 2026-02-15" "Generated: 2026-02-22"

Replace-Text "Total Federal Climate Resilience Investment: `$259,216,178 across 10 Tribal Nations" `
             "Total Federal Climate Resilience Investment: `$235,915,977 across 11 Tribal Nations"

Replace-Text "Aggregate Economic Impact: `$571,475,120 to `$761,966,826" `
             "Aggregate Economic Impact: `$528,589,758 to `$704,786,344"

Replace-Text "Estimated Jobs Supported: 2,540 to 4,762" "Estimated Jobs Supported: 2,349 to 4,405"

Replace-Text "Total Federal Climate Resilience Awards: `$259,216,178" `
             "Total Federal Climate Resilience Awards: `$235,915,977"

Replace-Text "Tribal Nations with Awards: 10 of 19 (53%)" "Tribal Nations with Awards: 11 of 19 (58%)"

Replace-Text "Investment Gap: 9 Tribal Nation(s) in this region have received zero federal climate resilience funding through tracked programs." `
             "Investment Gap: 8 Tribal Nation(s) in this region have received zero federal climate resilience funding through tracked programs."
